# Fruta / hortaliza, semanal
#
# Insert two new weekly report rows into the "Papa" (potato) price sheet,
# right before the existing row 358 (i.e. at position 356), shifting the
# rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 356 downward (existing rows 356..390 become 358..392).
$ws.Rows("356:357").Insert()

# --- New row 356 ---
$ws.Cells.Item(356, 1).Value = 5
$ws.Cells.Item(356, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(356, 3).Value = "Maule"
$ws.Cells.Item(356, 4).Value = 44578
$ws.Cells.Item(356, 5).Value = 7
$ws.Cells.Item(356, 6).Value = 100114001
$ws.Cells.Item(356, 7).Value = "Papa"
$ws.Cells.Item(356, 8).Value = "Asterix"
$ws.Cells.Item(356, 9).Value = "1a nueva(o)"
$ws.Cells.Item(356, 10).Value = 600
$ws.Cells.Item(356, 11).Value = 4000
$ws.Cells.Item(356, 12).Value = 4000
$ws.Cells.Item(356, 13).Value = 4000
$ws.Cells.Item(356, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(356, 15).Value = "Región del Maule"
$ws.Cells.Item(356, 16).Value = 160
$ws.Cells.Item(356, 17).Value = 25
$ws.Cells.Item(356, 18).Value = "Hortaliza"

# --- New row 357 ---
$ws.Cells.Item(357, 1).Value = 5
$ws.Cells.Item(357, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(357, 3).Value = "Maule"
$ws.Cells.Item(357, 4).Value = 44578
$ws.Cells.Item(357, 5).Value = 7
$ws.Cells.Item(357, 6).Value = 100114001
$ws.Cells.Item(357, 7).Value = "Papa"
$ws.Cells.Item(357, 8).Value = "Rosara"
$ws.Cells.Item(357, 9).Value = "1a nueva(o)"
$ws.Cells.Item(357, 10).Value = 800
$ws.Cells.Item(357, 11).Value = 7000
$ws.Cells.Item(357, 12).Value = 7000
$ws.Cells.Item(357, 13).Value = 7000
$ws.Cells.Item(357, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(357, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(357, 16).Value = 280
$ws.Cells.Item(357, 17).Value = 25
$ws.Cells.Item(357, 18).Value = "Hortaliza"
